$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI Gnai2->Adcy1 edge table: Sending cluster x Target cluster (FAPs, sCs)
# following Dr Hou advice -- recomputed with 3 target clusters considered per pair
$rows = @(
    @{ Row=2; A="ECs"; D="FAPs"; E=3; F=1; G=106.8785896666667; H=320.635769; I=0.1508748302900445; J=0.1508748302900445; K=3; L=1; M=0.1212106666666667; N=0.363632; O=0.774553595201428; P=0.774553595201428; Q=12.95482510588978; R=116.593425953008; S=0.1168606422265593; T=0.1168606422265593 },
    @{ Row=3; A="ECs"; D="sCs"; E=3; F=1; G=106.8785896666667; H=320.635769; I=0.1508748302900445; J=0.1508748302900445; K=2; L=0.6666666666666666; M=0.03528033333333334; N=0.105841; O=0.225446404798572; P=0.225446404798572; Q=3.770712269636555; R=33.936410426729; S=0.03401418806348523; T=0.03401418806348523 },
    @{ Row=4; A="FAPs"; D="FAPs"; E=3; F=1; G=77.232923; H=231.698769; I=0.1090256166999485; J=0.1090256166999485; K=3; L=1; M=0.1212106666666667; N=0.363632; O=0.774553595201428; P=0.774553595201428; Q=9.361454085445335; R=84.25308676900801; S=0.08444618338399797; T=0.08444618338399797 },
    @{ Row=5; A="FAPs"; D="sCs"; E=3; F=1; G=77.232923; H=231.698769; I=0.1090256166999485; J=0.1090256166999485; K=2; L=0.6666666666666666; M=0.03528033333333334; N=0.105841; O=0.225446404798572; P=0.225446404798572; Q=2.724803267747667; R=24.523229409729; S=0.02457943331595054; T=0.02457943331595054 },
    @{ Row=6; A="M1"; D="FAPs"; E=3; F=1; G=136.676337; H=410.029011; I=0.1929387280825172; J=0.1929387280825172; K=3; L=1; M=0.1212106666666667; N=0.363632; O=0.774553595201428; P=0.774553595201428; Q=16.566629925328; R=149.099669327952; S=0.1494413854899044; T=0.1494413854899044 },
    @{ Row=7; A="M1"; D="sCs"; E=3; F=1; G=136.676337; H=410.029011; I=0.1929387280825172; J=0.1929387280825172; K=2; L=0.6666666666666666; M=0.03528033333333334; N=0.105841; O=0.225446404798572; P=0.225446404798572; Q=4.821986728139001; R=43.39788055325101; S=0.04349734259261279; T=0.04349734259261279 },
    @{ Row=8; A="M2"; D="FAPs"; E=3; F=1; G=143.4557186666667; H=430.367156; I=0.2025088212285795; J=0.2025088212285795; K=3; L=1; M=0.1212106666666667; N=0.363632; O=0.774553595201428; P=0.774553595201428; Q=17.38836329673244; R=156.495269670592; S=0.1568539355425995; T=0.1568539355425995 },
    @{ Row=9; A="M2"; D="sCs"; E=3; F=1; G=143.4557186666667; H=430.367156; I=0.2025088212285795; J=0.2025088212285795; K=2; L=0.6666666666666666; M=0.03528033333333334; N=0.105841; O=0.225446404798572; P=0.225446404798572; Q=5.061165573132889; R=45.550490158196; S=0.04565488568597999; T=0.04565488568597999 },
    @{ Row=10; A="Neutro"; D="FAPs"; E=3; F=1; G=198.5982106666667; H=595.794632; I=0.2803505493821544; J=0.2803505493821544; K=3; L=1; M=0.1212106666666667; N=0.363632; O=0.774553595201428; P=0.774553595201428; Q=24.07222151371378; R=216.649993623424; S=0.2171465259406431; T=0.2171465259406431 },
    @{ Row=11; A="Neutro"; D="sCs"; E=3; F=1; G=198.5982106666667; H=595.794632; I=0.2803505493821544; J=0.2803505493821544; K=2; L=0.6666666666666666; M=0.03528033333333334; N=0.105841; O=0.225446404798572; P=0.225446404798572; Q=7.006611071723556; R=63.059499645512; S=0.06320402344151123; T=0.06320402344151123 },
    @{ Row=12; A="sCs"; D="FAPs"; E=3; F=1; G=45.55066433333334; H=136.651993; I=0.06430145431675577; J=0.06430145431675577; K=3; L=1; M=0.1212106666666667; N=0.363632; O=0.774553595201428; P=0.774553595201428; Q=5.52122639095289; R=49.691037518576; S=0.04980492261772356; T=0.04980492261772356 },
    @{ Row=13; A="sCs"; D="sCs"; E=3; F=1; G=45.55066433333334; H=136.651993; I=0.06430145431675577; J=0.06430145431675577; K=2; L=0.6666666666666666; M=0.03528033333333334; N=0.105841; O=0.225446404798572; P=0.225446404798572; Q=1.607042621234778; R=14.463383591113; S=0.01449653169903221; T=0.01449653169903221 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = "Gnai2"
    $ws.Range("C$row").Value = "Adcy1"
    $ws.Range("D$row").Value = $r.D
    foreach ($col in @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")) {
        $ws.Range("$col$row").Value = $r[$col]
    }
}